$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at 22 (old rows 22+ shift down by one: old 26/27 -> 27/28)
$ws.Rows.Item(22).Insert()

# 2) Copy the (old) last data row (row 21, still holding the 2508/ALEXANDER record with
#    the special "closing" border style) down into the newly inserted row 22 - this carries
#    both the values and the "last row" formatting to row 22.
$ws.Range("B21:J21").Copy($ws.Range("B22:J22"))

# 3) Row 21 is no longer the last row of the table, so it should pick up the regular
#    "interior" row formatting (same as rows 16-20) instead of the special closing style.
#    Copy formats only from row 20 onto row 21 (values in row21 are left untouched).
$ws.Range("B20:J20").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 4) Fill in the new row 22 with the new worker record (period 2509)
$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1063160955"
$ws.Range("D22").Value = "ALEXANDER MANUEL RIVAS MORELO"
$ws.Range("E22").Value = "2509"
$ws.Range("F22").Value = 56940
$ws.Range("G22").Value = 1423500
$ws.Range("H22").Value = ""
$ws.Range("I22").Value = ""
$ws.Range("J22").Value = ""

# 5) Update the summary totals
$ws.Range("E11").Value = 349813
$ws.Range("F13").Value = 7
